$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Workbook window size (workbook.xml bookViews)
# ---------------------------------------------------------------------------
$excel.Width = 28800
$excel.Height = 17500

# ---------------------------------------------------------------------------
# 2. Shared string fix: "GPTScore" (K70) -> cleared; new header labels added
#    throughout the new rows below. We clear K70's old text but keep format.
# ---------------------------------------------------------------------------
$ws.Range("K70").ClearContents()

# ---------------------------------------------------------------------------
# 3. Column widths
#    Engine quantizes ColumnWidth to 1/6-character steps, so we feed the
#    nearest reachable input for each target stored width:
#      col B  (2) -> stored 19.33203125  (closest reachable 19.333333...)
#      col I  (9) -> stored 15.83203125  (closest reachable 15.833333...)
#      col J (10) -> stored 14.5          (exact)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 15.0
$ws.Columns.Item(10).ColumnWidth = 13.666666666666666

# ---------------------------------------------------------------------------
# 4. Build the two new small-font styles used throughout the new block
#    (fontId 3 = sz8 Calibri, no border):
#      styleSmall       -> plain, used for header/avg numeric cells
#      styleSmallRight  -> right aligned, used for the "avg" row label
#    We materialise them once on scratch cells, then Copy/PasteSpecial
#    (Formats) from those scratch cells everywhere we need them, so every
#    target cell shares the very same style index (no style duplication).
# ---------------------------------------------------------------------------
$scratchPlain = $ws.Range("M1")
$scratchPlain.Value = "x"
$scratchPlain.Font.Size = 8

$scratchRight = $ws.Range("M2")
$scratchRight.Value = "x"
$scratchRight.Font.Size = 8
$scratchRight.HorizontalAlignment = -4152   # xlRight

function Set-PlainStyle($rng) {
    $scratchPlain.Copy() | Out-Null
    $rng.PasteSpecial(-4122) | Out-Null       # xlPasteFormats
}

function Set-RightStyle($rng) {
    $scratchRight.Copy() | Out-Null
    $rng.PasteSpecial(-4122) | Out-Null       # xlPasteFormats
}

# ---------------------------------------------------------------------------
# 5. New block #1 : rows 77-89 ("eval: gpt-4" / similarity & reasonable & valid)
# ---------------------------------------------------------------------------
foreach ($col in @("C","D","F","G","I","J")) { Set-PlainStyle $ws.Range($col + "77") }
foreach ($col in @("C","D","F","G","I","J")) { Set-PlainStyle $ws.Range($col + "78") }

$ws.Range("C77").Value = "eval: gpt-4"
$ws.Range("F77").Value = "eval: gpt-4"
$ws.Range("I77").Value = "eval: gpt-4"

$ws.Range("C78").Value = "similarity - 3.5"
$ws.Range("D78").Value = "similarity - 4"
$ws.Range("F78").Value = "reasonable - 3.5"
$ws.Range("G78").Value = "reasonable - 4"
$ws.Range("I78").Value = "valid - 3.5"
$ws.Range("J78").Value = "valid - 4"

$data77_88 = @{
    79 = @{C=0.25;   D=0.4;    F=0.21;   G=0.45;  I=0.175; J=0.4}
    80 = @{C=0.7;    D=0.6;    F=0.7;    G=0.65;  I=0.55;  J="N/A"}
    81 = @{C=0.475;  D=0.3;    F=0.74;   G=0.6;   I=0.6;   J=0.65}
    82 = @{C=0.65;   D=0.78;   F=0.75;   G=0.5;   I=0.6;   J="N/A"}
    83 = @{C=0.65;   D=0.65;   F=0.75;   G=0.68;  I=0.58;  J=0.7}
    84 = @{C=0.775;  D=0.8292; F=0.625;  G=0.65;  I=0.765; J=0.75}
    85 = @{C=0.22;   D=0.25;   F=0.45;   G=0.65;  I=0.6;   J=0.175}
    86 = @{C=0.63;   D=0.65;   F=0.3;    G=0.5;   I=0.67;  J=0.65}
    87 = @{C=0.7;    D=0.4;    F=0.675;  G=0.6;   I=0.575; J=0.5}
    88 = @{C=0.725;  D=0.75;   F=0.875;  G=0.725; I=0.825; J=0.8}
}
foreach ($r in 79..88) {
    $row = $data77_88[$r]
    foreach ($col in @("C","D","F","G","I","J")) {
        $ws.Range($col + $r).Value = $row[$col]
    }
}

Set-RightStyle $ws.Range("B89")
foreach ($col in @("C","D","F","G","I","J")) { Set-PlainStyle $ws.Range($col + "89") }

$ws.Range("B89").Value = "avg"
$ws.Range("C89").Formula = "=AVERAGE(C79:C88)"
$ws.Range("D89").Formula = "=AVERAGE(D79:D88)"
$ws.Range("F89").Formula = "=AVERAGE(F79:F88)"
$ws.Range("G89").Formula = "=AVERAGE(G79:G88)"
$ws.Range("I89").Formula = "=AVERAGE(I79:I88)"
$ws.Range("J89").Formula = "=AVERAGE(J79,J81,J83:J88)"

# ---------------------------------------------------------------------------
# 6. New block #2 : rows 92-104 (innovation)
# ---------------------------------------------------------------------------
foreach ($col in @("C","D","F","G")) { Set-PlainStyle $ws.Range($col + "92") }
foreach ($col in @("C","D","F","G")) { Set-PlainStyle $ws.Range($col + "93") }

$ws.Range("C92").Value = "eval: gpt-4"

$ws.Range("C93").Value = "innovation - 3.5"
$ws.Range("D93").Value = "innovation - 4"

$data94_103 = @{
    94  = @{C=0.25; D=0.27}
    95  = @{C=0.33; D=0.6}
    96  = @{C=0;    D=0.125}
    97  = @{C=0.1;  D=0.25}
    98  = @{C=0.6;  D=0.15}
    99  = @{C=0.3;  D=0.4}
    100 = @{C=0.44; D=0.7}
    101 = @{C=0.5;  D=0.6}
    102 = @{C=0.4;  D=0.4}
    103 = @{C=0.15; D=0.3}
}
foreach ($r in 94..103) {
    $row = $data94_103[$r]
    $ws.Range("C" + $r).Value = $row["C"]
    $ws.Range("D" + $r).Value = $row["D"]
}

Set-RightStyle $ws.Range("B104")
Set-RightStyle $ws.Range("E104")
foreach ($col in @("C","D","F","G")) { Set-PlainStyle $ws.Range($col + "104") }

$ws.Range("B104").Value = "avg"
$ws.Range("C104").Formula = "=AVERAGE(C94:C103)"
$ws.Range("D104").Formula = "=AVERAGE(D94:D103)"

# ---------------------------------------------------------------------------
# 7. Clean up scratch cells used to build the reusable styles
# ---------------------------------------------------------------------------
$ws.Range("M1").Clear() | Out-Null
$ws.Range("M2").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 8. View state: scroll position, zoom, selection (best effort)
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 144
$excel.ActiveWindow.ScrollRow = 72
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F92").Select()
